$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol data (price + 1h volume change) scraped on 2023-02-08.
# Values are kept as text (NumberFormat "@") so exact string formatting
# (trailing zeros, decimal precision, "%" suffix) is preserved verbatim.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.69%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.35%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.262"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.69%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08336"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.91%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.942"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.68%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9717"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.13%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1895"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.34%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09653"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.13%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04621"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.32%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1060"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.51%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001283"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006022"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.71%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.402"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.84%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.455"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.65%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3362"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.69%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.661"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-14.88%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1364"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.94%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2582"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.43%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04156"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.20%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.23%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004427"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.84%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.93%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002988"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.03%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02710"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05656"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.34%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007837"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1411"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007362"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.98%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007874"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.13%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3496"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006879"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.13%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.35%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003494"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.44%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003541"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.64%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.35%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.35%"
